$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Reference cell with the default (unstyled) cell style, used to reset
# style after coercing ambiguous numeric-looking strings to stay as text
# (matches the source file, which stores these as plain text with no style).
$defaultStyle = $ws.Range("B2").Style

# Row 2
$ws.Range("D2").Value = "98.938.57"
$ws.Range("E2").Value = "  +1.06%  "

# Row 3
$ws.Range("D3").Value = "3.292.42"
$ws.Range("E3").Value = "  -1.65%  "

# Row 4
$ws.Range("E4").Value = "  -0.07%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "254.71"
$ws.Range("D5").Style = $defaultStyle
$ws.Range("E5").Value = "  +0.29%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "627.29"
$ws.Range("D6").Style = $defaultStyle
$ws.Range("E6").Value = "  +0.95%  "

# Row 7
$ws.Range("E7").Value = "  +23.07%  "

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.402"
$ws.Range("D8").Style = $defaultStyle
$ws.Range("E8").Value = "  +4.74%  "

# Row 9
$ws.Range("E9").Value = "  -0.06%  "

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.985"
$ws.Range("D10").Style = $defaultStyle
$ws.Range("E10").Value = "  +23.48%  "

# Row 11
$ws.Range("D11").Value = "3.289.67"
$ws.Range("E11").Value = "  -1.58%  "

# Row 12
$ws.Range("E12").Value = "  +2.73%  "

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "40.83"
$ws.Range("D13").Style = $defaultStyle
$ws.Range("E13").Value = "  +14.44%  "

# Row 14
$ws.Range("D14").Value = "98.718.44"
$ws.Range("E14").Value = "  +1.04%  "

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.0000250"
$ws.Range("D15").Style = $defaultStyle
$ws.Range("E15").Value = "  +2.04%  "

# Row 16
$ws.Range("D16").Value = "3.908.23"
$ws.Range("E16").Value = "  -1.64%  "

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "5.48"
$ws.Range("D17").Style = $defaultStyle
$ws.Range("E17").Value = "  +0.21%  "

# Row 18
$ws.Range("D18").Value = "3.293.16"
$ws.Range("E18").Value = "  -1.47%  "

# Row 19
$ws.Range("E19").Value = "  -3.48%  "

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "15.57"
$ws.Range("D20").Style = $defaultStyle
$ws.Range("E20").Value = "  +5.76%  "

# Row 21
$ws.Range("E21").Value = "  +9.28%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "488.32"
$ws.Range("D22").Style = $defaultStyle
$ws.Range("E22").Value = "  +2.04%  "

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "9.40"
$ws.Range("D23").Style = $defaultStyle
$ws.Range("E23").Value = "  +3.49%  "

# Row 24
$ws.Range("E24").Value = "  -1.75%  "

# Row 25
$ws.Range("B25").Value = "NEARProtocol"
$ws.Range("C25").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "5.71"
$ws.Range("D25").Style = $defaultStyle
$ws.Range("E25").Value = "  +0.81%  "

# Row 26
$ws.Range("B26").Value = "Stellar"
$ws.Range("C26").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.343"
$ws.Range("D26").Style = $defaultStyle
$ws.Range("E26").Value = "  +39.32%  "

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "89.02"
$ws.Range("D27").Style = $defaultStyle
$ws.Range("E27").Value = "  +1.72%  "

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "12.16"
$ws.Range("D28").Style = $defaultStyle
$ws.Range("E28").Value = "  +2.09%  "

# Row 29
$ws.Range("D29").Value = "3.462.99"
$ws.Range("E29").Value = "  -3.03%  "

# Row 30
$ws.Range("E30").Value = "  +19.07%  "

# Row 31
$ws.Range("E31").Value = "  -0.09%  "

# Row 32
$ws.Range("E32").Value = "  +1.10%  "

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "10.73"
$ws.Range("D33").Style = $defaultStyle
$ws.Range("E33").Value = "  +17.68%  "

# Row 34
$ws.Range("E34").Value = "  -0.14%  "

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "27.95"
$ws.Range("D35").Style = $defaultStyle
$ws.Range("E35").Value = "  +3.01%  "

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.481"
$ws.Range("D36").Style = $defaultStyle
$ws.Range("E36").Value = "  +7.69%  "

# Row 37
$ws.Range("E37").Value = "  +0.46%  "

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "7.32"
$ws.Range("D38").Style = $defaultStyle
$ws.Range("E38").Value = "  +1.09%  "

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.96"
$ws.Range("D39").Style = $defaultStyle
$ws.Range("E39").Value = "  +1.19%  "

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "495.66"
$ws.Range("D40").Style = $defaultStyle
$ws.Range("E40").Value = "  -4.67%  "

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "24.73"
$ws.Range("D41").Style = $defaultStyle

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "3.83"
$ws.Range("D42").Style = $defaultStyle
$ws.Range("E42").Value = "  +5.77%  "

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.25"
$ws.Range("D43").Style = $defaultStyle
$ws.Range("E43").Value = "  -1.11%  "

# Row 44
$ws.Range("E44").Value = "  +0.01%  "

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.783"
$ws.Range("D45").Style = $defaultStyle
$ws.Range("E45").Value = "  -0.80%  "

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "3.15"
$ws.Range("D46").Style = $defaultStyle
$ws.Range("E46").Value = "  -2.06%  "

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "158.91"
$ws.Range("D47").Style = $defaultStyle
$ws.Range("E47").Value = "  -1.21%  "

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.95"
$ws.Range("D48").Style = $defaultStyle
$ws.Range("E48").Value = "  +2.28%  "

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "4.83"
$ws.Range("D49").Style = $defaultStyle
$ws.Range("E49").Value = "  +7.76%  "

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "7.36"
$ws.Range("D50").Style = $defaultStyle
$ws.Range("E50").Value = "  +16.42%  "

# Row 51
$ws.Range("E51").Value = "  +7.67%  "
